# Generate Report for Handback
# Both handback files (af32eaaf... and f8ddade8...) have been handed back;
# update Status columns from "Ready for handoff" to
# "Handed back: in sync with en-US" across all three sheets, and refresh the
# Latest Handback DateTime / Error Detail for the f8ddade8 file on the
# zh-cn and de-de report sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
# Row 3 corresponds to f8ddade8-8c09-481e-95e5-1efd3d21c69f.md, whose
# zh-cn/de-de status columns still said "Ready for handoff".
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet -------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
# Row 2: af32eaaf-086f-49c9-b5ae-da34c1160d0e.md
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
# Row 3: f8ddade8-8c09-481e-95e5-1efd3d21c69f.md
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("K3").Value = "2016-09-02 04:52:01"
$zhcn.Range("P3").Value = ""

# --- de-de sheet -------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
# Row 3: f8ddade8-8c09-481e-95e5-1efd3d21c69f.md
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("K3").Value = "2016-09-02 04:52:17"
$dede.Range("P3").Value = ""
